$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 46045
$ws.Cells.Item(2, 2).Value = 12762.036064577
$ws.Cells.Item(2, 3).Value = 11703.0628761372
$ws.Cells.Item(2, 4).Value = 18147.86
$ws.Cells.Item(2, 5).Value = 8325.91277327107
$ws.Cells.Item(2, 6).Value = 78.3798187253451

$ws.Cells.Item(3, 1).Value = 46046
$ws.Cells.Item(3, 2).Value = 4899.73315221571
$ws.Cells.Item(3, 3).Value = 7729.43829099819
$ws.Cells.Item(3, 4).Value = 12075.86
$ws.Cells.Item(3, 5).Value = 8009.55992026545
$ws.Cells.Item(3, 6).Value = 152.630758802652

$ws.Cells.Item(4, 1).Value = 46047
$ws.Cells.Item(4, 2).Value = 4688.14215797149
$ws.Cells.Item(4, 3).Value = 7631.97900650917
$ws.Cells.Item(4, 4).Value = 12075.86
$ws.Cells.Item(4, 5).Value = 7866.01334327686
$ws.Cells.Item(4, 6).Value = 142.588847907751

$ws.Cells.Item(5, 1).Value = 46048
$ws.Cells.Item(5, 2).Value = 12095.4314989621
$ws.Cells.Item(5, 3).Value = 11729.0039029181
$ws.Cells.Item(5, 4).Value = 12075.86
$ws.Cells.Item(5, 5).Value = 8119.30813408507
$ws.Cells.Item(5, 6).Value = 323.852168208466

$ws.Cells.Item(6, 1).Value = 46049
$ws.Cells.Item(6, 2).Value = 11870.5678861241
$ws.Cells.Item(6, 3).Value = 11934.2928006436
$ws.Cells.Item(6, 4).Value = 12075.86
$ws.Cells.Item(6, 5).Value = 7961.24071870759
$ws.Cells.Item(6, 6).Value = 325.819729972966

$ws.Cells.Item(7, 1).Value = 46050
$ws.Cells.Item(7, 2).Value = 12315.9682835607
$ws.Cells.Item(7, 3).Value = 11593.8793174649
$ws.Cells.Item(7, 4).Value = 12075.86
$ws.Cells.Item(7, 5).Value = 8314.35920730669
$ws.Cells.Item(7, 6).Value = 326.349105198818

$ws.Cells.Item(8, 1).Value = 46051
$ws.Cells.Item(8, 2).Value = 12315.9682835607
$ws.Cells.Item(8, 3).Value = 11832.167532391
$ws.Cells.Item(8, 4).Value = 12075.86
$ws.Cells.Item(8, 5).Value = 8314.35920730669
$ws.Cells.Item(8, 6).Value = 336.277780820737

$ws.Cells.Item(9, 1).Value = 46052
$ws.Cells.Item(9, 2).Value = 12315.9682835607
$ws.Cells.Item(9, 3).Value = 11395.7558317288
$ws.Cells.Item(9, 4).Value = 12075.86
$ws.Cells.Item(9, 5).Value = 8314.35920730669
$ws.Cells.Item(9, 6).Value = 318.093959959811

$ws.Cells.Item(10, 1).Value = 46053
$ws.Cells.Item(10, 2).Value = 4867.38022112383
$ws.Cells.Item(10, 3).Value = 7991.94528864509
$ws.Cells.Item(10, 4).Value = 12075.86
$ws.Cells.Item(10, 5).Value = 7930.29742224961
$ws.Cells.Item(10, 6).Value = 160.265946287279

$ws.Cells.Item(11, 1).Value = 46054
$ws.Cells.Item(11, 2).Value = 5046.61366744637
$ws.Cells.Item(11, 3).Value = 7842.65473127017
$ws.Cells.Item(11, 4).Value = 9743.86
$ws.Cells.Item(11, 5).Value = 7892.99600023043
$ws.Cells.Item(11, 6).Value = 249.657947145858

$ws.Cells.Item(12, 1).Value = 46055
$ws.Cells.Item(12, 2).Value = 11445.6138712783
$ws.Cells.Item(12, 3).Value = 11220.8700408412
$ws.Cells.Item(12, 4).Value = 9743.86
$ws.Cells.Item(12, 5).Value = 7808.58662750625
$ws.Cells.Item(12, 6).Value = 386.899861181144

$ws.Cells.Item(13, 1).Value = 46056
$ws.Cells.Item(13, 2).Value = 11445.6138712783
$ws.Cells.Item(13, 3).Value = 11417.6235115178
$ws.Cells.Item(13, 4).Value = 9743.86
$ws.Cells.Item(13, 5).Value = 7808.58662750625
$ws.Cells.Item(13, 6).Value = 395.097922459336

$ws.Cells.Item(14, 1).Value = 46057
$ws.Cells.Item(14, 2).Value = 11445.6138712783
$ws.Cells.Item(14, 3).Value = 11248.4911458367
$ws.Cells.Item(14, 4).Value = 9743.86
$ws.Cells.Item(14, 5).Value = 7808.58662750625
$ws.Cells.Item(14, 6).Value = 388.050740555958

$ws.Cells.Item(15, 1).Value = 46058
$ws.Cells.Item(15, 2).Value = 11445.6138712783
$ws.Cells.Item(15, 3).Value = 10936.8445968102
$ws.Cells.Item(15, 4).Value = 9743.86
$ws.Cells.Item(15, 5).Value = 7808.58662750625
$ws.Cells.Item(15, 6).Value = 375.065467679852
